# Apply the "Add files via upload" update to the Saldo export sheet.
#
# Summary of the change:
#  - Three account balances were updated (MERG, CASSIO, MARIA).
#  - The sheet is kept sorted by Saldo (balance) descending, so after the
#    balance updates the first block of data rows (accounts 004212581
#    .. 004461070, currently rows 2-13) needs to be re-sorted.
#  - A later row (JOAO / 004381328, row 15) also has its balance updated,
#    but it stays in place because the new value still fits between its
#    neighbours in the descending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the balances that changed (value-only edits keep the existing
# cell type/formatting for the Conta/Nome text columns intact).
$ws.Cells.Item(4, 3).Value  = 236549.41   # MERG    (004214592)
$ws.Cells.Item(10, 3).Value = 73331.58    # CASSIO  (004508526)
$ws.Cells.Item(2, 3).Value  = 67263.08    # MARIA   (004212581)
$ws.Cells.Item(15, 3).Value = 52000       # JOAO    (004381328)

# Re-sort the top block (rows 2-13) descending by Saldo (column C) so the
# updated balances land in the correct position again.
$rng = $ws.Range("A2:C13")
$rng.Sort($ws.Range("C2:C13"), 2)
